$p = $ppt.ActivePresentation

# --- Slide 2 ("Creation of SASS") is being repurposed into a new "Goal" slide.
# First, duplicate it so the duplicate (inserted right after, at index 3) keeps
# the original "Creation of SASS" title + bullet content.
$original = $p.Slides.Item(2)
[void]$original.Duplicate()

# --- Now update the original slide (still at index 2) with the new "Goal" content.
$goalSlide = $p.Slides.Item(2)

$titleShape = $goalSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Goal"

$bodyShape = $goalSlide.Shapes.Item(2)
$bodyShape.TextFrame.TextRange.Text = "Introduce SASS and its features and benefits that it brings to CSS`rDiscuss how it is compiled to CSS for the browser "
[void]$bodyShape.TextFrame.TextRange.InsertAfter("to understand")
